$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 14:05"

# Country reorderings + refreshed case counts
# Row 26
$ws.Range("B26").Value = 32426
$ws.Range("C26").Value = 918
$ws.Range("D26").Value = 11415
$ws.Range("E26").Value = 20832
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 179

# Row 37
$ws.Range("A37").Value = "Kuwait"
$ws.Range("B37").Value = 17568
$ws.Range("C37").Value = 804
$ws.Range("D37").Value = 4885
$ws.Range("E37").Value = 12559
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 124

# Row 38
$ws.Range("A38").Value = "Rumania"
$ws.Range("B38").Value = 17387
$ws.Range("C38").Value = 196
$ws.Range("D38").Value = 10356
$ws.Range("E38").Value = 5890
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 1141

# Row 39
$ws.Range("A39").Value = "Sudafrica"
$ws.Range("B39").Value = 17200
$ws.Range("D39").Value = 7960
$ws.Range("E39").Value = 8928
$ws.Range("H39").Value = 312

# Row 40
$ws.Range("A40").Value = "Colombia"
$ws.Range("B40").Value = 16935
$ws.Range("D40").Value = 4050
$ws.Range("E40").Value = 12272
$ws.Range("H40").Value = 613

# Row 55
$ws.Range("B55").Value = 7843
$ws.Range("C55").Value = 311
$ws.Range("D55").Value = 3334
$ws.Range("E55").Value = 4497

# Row 186
$ws.Range("A186").Value = "Comoras"
$ws.Range("B186").Value = 34
$ws.Range("C186").Value = 23
$ws.Range("D186").Value = 8
$ws.Range("E186").Value = 25

# Row 187
$ws.Range("A187").Value = "Guam"
$ws.Range("B187").Value = 32
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 31

# Row 188
$ws.Range("A188").Value = "Botsuana"
$ws.Range("D188").Value = 17
$ws.Range("E188").Value = 7
$ws.Range("H188").Value = 1

# Row 189
$ws.Range("A189").Value = "Antigua y Barbuda"
$ws.Range("B189").Value = 25
$ws.Range("D189").Value = 19
$ws.Range("E189").Value = 3
$ws.Range("H189").Value = 3

# Row 190
$ws.Range("A190").Value = "Gambia"
$ws.Range("D190").Value = 13
$ws.Range("E190").Value = 10
$ws.Range("H190").Value = 1

# Row 191
$ws.Range("A191").Value = "Timor Oriental"
$ws.Range("B191").Value = 24
$ws.Range("D191").Value = 24
$ws.Range("E191").Value = 0

# Row 192
$ws.Range("A192").Value = "Granada"
$ws.Range("B192").Value = 22
$ws.Range("D192").Value = 14
$ws.Range("E192").Value = 8

# Row 193
$ws.Range("A193").Value = "Butan"
$ws.Range("B193").Value = 21
$ws.Range("D193").Value = 5
$ws.Range("E193").Value = 16

# Row 194
$ws.Range("A194").Value = "Laos"
$ws.Range("B194").Value = 19
$ws.Range("D194").Value = 14
$ws.Range("E194").Value = 5

# Row 195
$ws.Range("A195").Value = "Fiyi"
$ws.Range("D195").Value = 15
$ws.Range("E195").Value = 3
$ws.Range("H195").Value = 0

# Row 197
$ws.Range("A197").Value = "Belice"
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2

# Row 198
$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("B198").Value = 18
$ws.Range("D198").Value = 18
$ws.Range("E198").Value = 0

# Row 199
$ws.Range("A199").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 17

# Row 200
$ws.Range("A200").Value = "San Vicente y las Granadinas"
$ws.Range("B200").Value = 17
$ws.Range("E200").Value = 3

# Row 201
$ws.Range("A201").Value = "Namibia"
$ws.Range("E201").Value = 2
$ws.Range("H201").Value = 0

# Row 202
$ws.Range("A202").Value = "Curazao"
$ws.Range("D202").Value = 14
$ws.Range("E202").Value = 1
$ws.Range("H202").Value = 1

# Row 203
$ws.Range("A203").Value = "Dominica"
$ws.Range("B203").Value = 16
$ws.Range("D203").Value = 16

# Row 204
$ws.Range("A204").Value = "San Cristobal y Nieves"
$ws.Range("B204").Value = 15
$ws.Range("D204").Value = 15

# Row 205
$ws.Range("A205").Value = "Islas Malvinas"
$ws.Range("B205").Value = 13
$ws.Range("D205").Value = 13
$ws.Range("E205").Value = 0

# Row 206
$ws.Range("A206").Value = "Santa Sede"
$ws.Range("D206").Value = 2
$ws.Range("E206").Value = 10
$ws.Range("H206").Value = 0

# Row 207
$ws.Range("A207").Value = "Islas Turcas y Caicos"
$ws.Range("B207").Value = 12
$ws.Range("D207").Value = 10
$ws.Range("E207").Value = 1

# Row 209
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

# Row 210
$ws.Range("A210").Value = "Groenlandia"

# Row 211
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

